$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (style matches existing header row cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(3, 5),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 9),
    @(5, 7),
    @(10, 10),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(6, 8),
    @(9, 9),
    @(7, 8),
    @(1, 3),
    @(2, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
